$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 324.97827
$ws.Range("I8").Value = 559.75
$ws.Range("J8").Value = 302.61905
$ws.Range("K8").Value = 1679.25
$ws.Range("L8").Value = 907.85715
$ws.Range("M8").Value = -1540.25
$ws.Range("N8").Value = -1185.85715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4966.25
$ws.Range("I19").Value = 3732.6667
$ws.Range("K19").Value = 3732.6667
$ws.Range("M19").Value = -3557.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 38891284
$ws.Range("I70").Value = 50002212
$ws.Range("J70").Value = 33335820
$ws.Range("K70").Value = 150006636
$ws.Range("L70").Value = 100007460
$ws.Range("M70").Value = -150006366
$ws.Range("N70").Value = -100008000

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 38891284
$ws.Range("I73").Value = 50002212
$ws.Range("J73").Value = 33335820
$ws.Range("K73").Value = 150006636
$ws.Range("L73").Value = 100007460
$ws.Range("M73").Value = -150005700
$ws.Range("N73").Value = -100009332

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 3290416
$ws.Range("I107").Value = 4032973
$ws.Range("J107").Value = 1948.8572
$ws.Range("K107").Value = 4032973
$ws.Range("L107").Value = 1948.8572
$ws.Range("M107").Value = -4031053
$ws.Range("N107").Value = -5788.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 3376
$ws.Range("I127").Value = 3595.875
$ws.Range("K127").Value = 10787.625
$ws.Range("M127").Value = -5827.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1352.2941
$ws.Range("I132").Value = 1382.125
$ws.Range("K132").Value = 4146.375
$ws.Range("M132").Value = -1616.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3848.9805
$ws.Range("J138").Value = 7108.4346
$ws.Range("L138").Value = 21325.3038
$ws.Range("N138").Value = -31605.3038

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 79561.125
$ws.Range("J139").Value = 79641.42999999999
$ws.Range("L139").Value = 79641.42999999999
$ws.Range("N139").Value = -89921.42999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2816.9546
$ws.Range("I141").Value = 2598.85
$ws.Range("J141").Value = 4998
$ws.Range("K141").Value = 7796.549999999999
$ws.Range("L141").Value = 14994
$ws.Range("M141").Value = -2616.549999999999
$ws.Range("N141").Value = -25354

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1249
$ws.Range("I16").Value = 748.3333
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 748.3333
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -461.3333
$ws.Range("N16").Value = -2574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8224.467000000001
$ws.Range("I45").Value = 2454.2856
$ws.Range("K45").Value = 2454.2856
$ws.Range("M45").Value = -2077.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6478.086
$ws.Range("I61").Value = 2725.9473
$ws.Range("J61").Value = 10933.75
$ws.Range("K61").Value = 2725.9473
$ws.Range("L61").Value = 10933.75
$ws.Range("M61").Value = -2513.9473
$ws.Range("N61").Value = -11357.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2136.75
$ws.Range("I63").Value = 2326
$ws.Range("K63").Value = 2326
$ws.Range("M63").Value = -1640

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2136.75
$ws.Range("I66").Value = 2326
$ws.Range("K66").Value = 11630
$ws.Range("M66").Value = -8198

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14588.54
$ws.Range("I74").Value = 21425.4
$ws.Range("J74").Value = 4333.25
$ws.Range("K74").Value = 21425.4
$ws.Range("L74").Value = 4333.25
$ws.Range("M74").Value = -20551.4
$ws.Range("N74").Value = -6081.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 14588.54
$ws.Range("I77").Value = 21425.4
$ws.Range("J77").Value = 4333.25
$ws.Range("K77").Value = 107127
$ws.Range("L77").Value = 21666.25
$ws.Range("M77").Value = -102759
$ws.Range("N77").Value = -30402.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 22223438
$ws.Range("I110").Value = 1162.6666
$ws.Range("K110").Value = 1162.6666
$ws.Range("M110").Value = 882.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 11296.131
$ws.Range("I122").Value = 13649.8125
$ws.Range("K122").Value = 40949.4375
$ws.Range("M122").Value = -38499.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4848.25
$ws.Range("I132").Value = 2746.2
$ws.Range("K132").Value = 8238.599999999999
$ws.Range("M132").Value = -5708.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6478.086
$ws.Range("I136").Value = 2725.9473
$ws.Range("J136").Value = 10933.75
$ws.Range("K136").Value = 8177.841899999999
$ws.Range("L136").Value = 32801.25
$ws.Range("M136").Value = -5627.841899999999
$ws.Range("N136").Value = -37901.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I64").Value = 111111780
$ws.Range("J64").Value = 1233.3334
$ws.Range("K64").Value = 111111780
$ws.Range("L64").Value = 1233.3334
$ws.Range("M64").Value = -111111555
$ws.Range("N64").Value = -1683.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I67").Value = 111111780
$ws.Range("J67").Value = 1233.3334
$ws.Range("K67").Value = 111111780
$ws.Range("L67").Value = 1233.3334
$ws.Range("M67").Value = -111111000
$ws.Range("N67").Value = -2793.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3229.1035
$ws.Range("I105").Value = 3006.8635
$ws.Range("K105").Value = 3006.8635
$ws.Range("M105").Value = -1259.8635

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 43271536
$ws.Range("I107").Value = 59211370
$ws.Range("J107").Value = 6285.2856
$ws.Range("K107").Value = 59211370
$ws.Range("L107").Value = 6285.2856
$ws.Range("M107").Value = -59209450
$ws.Range("N107").Value = -10125.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3818.4348
$ws.Range("I16").Value = 2148
$ws.Range("J16").Value = 5990
$ws.Range("K16").Value = 2148
$ws.Range("L16").Value = 5990
$ws.Range("M16").Value = -1861
$ws.Range("N16").Value = -6564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 40750
$ws.Range("J64").Value = 40750
$ws.Range("L64").Value = 40750
$ws.Range("N64").Value = -41246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 40750
$ws.Range("J67").Value = 40750
$ws.Range("L67").Value = 40750
$ws.Range("N67").Value = -42466

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6445.773
$ws.Range("I99").Value = 5487.1875
$ws.Range("K99").Value = 5487.1875
$ws.Range("M99").Value = -3989.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1697.1389
$ws.Range("I107").Value = 1146.5238
$ws.Range("K107").Value = 1146.5238
$ws.Range("M107").Value = 773.4762000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3818.4348
$ws.Range("I113").Value = 2148
$ws.Range("J113").Value = 5990
$ws.Range("K113").Value = 2148
$ws.Range("L113").Value = 5990
$ws.Range("M113").Value = 22
$ws.Range("N113").Value = -10330

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1553.0526
$ws.Range("I122").Value = 1042.1
$ws.Range("K122").Value = 3126.3
$ws.Range("M122").Value = -676.2999999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6445.773
$ws.Range("I126").Value = 5487.1875
$ws.Range("K126").Value = 16461.5625
$ws.Range("M126").Value = -13991.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 79634970
$ws.Range("I55").Value = 222222260
$ws.Range("J55").Value = 8341322
$ws.Range("K55").Value = 666666780
$ws.Range("L55").Value = 25023966
$ws.Range("M55").Value = -666666603
$ws.Range("N55").Value = -25024320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 54289.58
$ws.Range("J137").Value = 127215.625
$ws.Range("L137").Value = 381646.875
$ws.Range("N137").Value = -391846.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 60563.938
$ws.Range("J57").Value = 60563.938
$ws.Range("L57").Value = 60563.938
$ws.Range("N57").Value = -62203.938

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 581.8
$ws.Range("I97").Value = 581.8
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 581.8
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -85.79999999999995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8601.77
$ws.Range("I132").Value = 2102.875
$ws.Range("K132").Value = 6308.625
$ws.Range("M132").Value = -3778.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5015.3477
$ws.Range("I40").Value = 4103.6924
$ws.Range("K40").Value = 4103.6924
$ws.Range("M40").Value = -3967.6924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1255.125
$ws.Range("I82").Value = 816.1667
$ws.Range("K82").Value = 816.1667
$ws.Range("M82").Value = -455.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1255.125
$ws.Range("I85").Value = 816.1667
$ws.Range("K85").Value = 816.1667
$ws.Range("M85").Value = 431.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3978.2
$ws.Range("I100").Value = 2845.5
$ws.Range("J100").Value = 5272.7144
$ws.Range("K100").Value = 2845.5
$ws.Range("L100").Value = 5272.7144
$ws.Range("M100").Value = -2304.5
$ws.Range("N100").Value = -6354.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5993.0625
$ws.Range("I122").Value = 3876
$ws.Range("J122").Value = 8715
$ws.Range("K122").Value = 11628
$ws.Range("L122").Value = 26145
$ws.Range("M122").Value = -9178
$ws.Range("N122").Value = -31045

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12507204
$ws.Range("I132").Value = 26318324
$ws.Range("K132").Value = 78954972
$ws.Range("M132").Value = -78952442

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 55556940
$ws.Range("I107").Value = 1000.3333
$ws.Range("K107").Value = 3000.9999
$ws.Range("M107").Value = -1080.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("N118").Value = 0
